# Natmi following Dr Hou advice
# Re-derive the LR-pair table for Lgi4-Adam11 across the full 3x3 grid of
# clusters (ECs, FAPs, sCs) as sending/target cluster instead of only the
# two rows that existed before (FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs",  "Lgi4", "Adam11", "ECs",  1, 0.3333333333333333, 0.03034166666666667, 0.09102499999999999, 0.002265018822860075, 0.002265018822860075, 2, 0.6666666666666666, 0.114041,            0.342123,            0.07100488927427132, 0.07100488927427133, 0.003460194008333333, 0.031141746075,     0.00016082741072132,  0.00016082741072132),
    @("ECs",  "Lgi4", "Adam11", "FAPs", 1, 0.3333333333333333, 0.03034166666666667, 0.09102499999999999, 0.002265018822860075, 0.002265018822860075, 2, 0.6666666666666666, 0.292138,            0.876414,            0.181892708261126,   0.1818927082611261,  0.008863953816666667, 0.07977558435,      0.0004119904079524468, 0.0004119904079524468),
    @("ECs",  "Lgi4", "Adam11", "sCs",  1, 0.3333333333333333, 0.03034166666666667, 0.09102499999999999, 0.002265018822860075, 0.002265018822860075, 3, 1,                   1.199921666666667,   3.599765,            0.7471024024646026,  0.7471024024646027,  0.03640762323611111,  0.327668609125,     0.001692201004186308,  0.001692201004186309),
    @("FAPs", "Lgi4", "Adam11", "ECs",  3, 1,                   10.326286,           30.978858,           0.7708618124768957,   0.7708618124768957,   2, 0.6666666666666666, 0.114041,            0.342123,            0.07100488927427132, 0.07100488927427133, 1.177619981726,       10.598579835534,    0.05473495764068607,   0.05473495764068609),
    @("FAPs", "Lgi4", "Adam11", "FAPs", 3, 1,                   10.326286,           30.978858,           0.7708618124768957,   0.7708618124768957,   2, 0.6666666666666666, 0.292138,            0.876414,            0.181892708261126,   0.1818927082611261,  3.016700539468,       27.150304855212,    0.1402141427665028,    0.1402141427665028),
    @("FAPs", "Lgi4", "Adam11", "sCs",  3, 1,                   10.326286,           30.978858,           0.7708618124768957,   0.7708618124768957,   3, 1,                   1.199921666666667,   3.599765,            0.7471024024646026,  0.7471024024646027,  12.39073430759666,    111.51660876837,    0.5759127120697067,    0.5759127120697068),
    @("sCs",  "Lgi4", "Adam11", "ECs",  3, 1,                   3.039140333333334,  9.117421,            0.2268731687002442,   0.2268731687002443,   2, 0.6666666666666666, 0.114041,            0.342123,            0.07100488927427132, 0.07100488927427133, 0.3465866027536667,   3.119279424783,     0.01610910422286392,   0.01610910422286392),
    @("sCs",  "Lgi4", "Adam11", "FAPs", 3, 1,                   3.039140333333334,  9.117421,            0.2268731687002442,   0.2268731687002443,   2, 0.6666666666666666, 0.292138,            0.876414,            0.181892708261126,   0.1818927082611261,  0.8878483786993334,   7.990635408294001,  0.04126657508667075,   0.04126657508667077),
    @("sCs",  "Lgi4", "Adam11", "sCs",  3, 1,                   3.039140333333334,  9.117421,            0.2268731687002442,   0.2268731687002443,   3, 1,                   1.199921666666667,   3.599765,            0.7471024024646026,  0.7471024024646027,  3.646730334007222,    32.820573006065,    0.1694974893907095,    0.1694974893907096)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowValues = $rows[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowValues[$c]
    }
}
